$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("P2") "320018569086"
Set-TextValue $ws.Range("Q2") "`$19.04"
Set-TextValue $ws.Range("R2") "PASS"
Set-TextValue $ws.Range("P3") "320018590118"
Set-TextValue $ws.Range("P5") "320018567576"
Set-TextValue $ws.Range("Q5") "`$43.07"
